# Fruta / hortaliza, semanal
#
# The underlying source data got re-synced and several rows' weekly price
# records shuffled position (same column layout, values now belong to a
# different row). Concretely, rows {2,5,7}, rows {3,4,12} and rows
# {8,10,11} each form a 3-way rotation of the columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg) and Q (Kg o Unidades).
#
# Capture every source row's values FIRST (into variables), then write
# them all out to their destination row, so the 3-way rotations don't
# clobber a value before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row it should receive its values from
$rowMap = @{2 = 7; 3 = 4; 4 = 12; 5 = 2; 7 = 5; 8 = 11; 10 = 8; 11 = 10; 12 = 3}
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values of every row referenced above.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$srcRow").Value2()
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Write each destination row from its recorded source-row snapshot.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowVals[$col]
    }
}

"Reassigned rows: $($rowMap.Keys -join ', ')"
